# Kandidat_2022_vars.xlsx - update question wording, drop an obsolete
# "heltidsstilling" row, tidy up a couple of comments / svartype values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Edits that happen on rows still above the row we are about to delete ---

# A2: shorten/rephrase the "hovedaktivitet" question text
$ws.Range("A2").Value = "Hovedaktivitet 1–3 år etter fullført utdanning"

# --- Remove the now-unused "Andel heltidsstilling" / "heltidsstilling" row ---
# (this is row 10 in the original layout; everything below shifts up by one)
$ws.Rows("10").Delete()

# --- Edits on rows below the deleted one, using their NEW (post-delete) row numbers ---

# E10 ("Ufrivillig deltid" comment): add an explanatory lead-in sentence
$ws.Range("E10").Value = "Andel av de som arbeider deltid og som ønsker full stilling. I 2018 fikk ikke kandidatene oppfølgingsspørsmål om grunnen til redusert stilling, og vi har derfor ikke tatt med sammenligning av resultater."

# C14 ("lang_tid_til_relevant_arbeid" svartype): now a regular (non-single) mean
$ws.Range("C14").Value = "snitt_as_num"

# A15:A16 - "fornoyd_oppgaver" question text, trimmed of the survey-statement wrapper
$ws.Range("A15").Value = "Jeg er godt fornøyd med oppgaver og ansvar i min nåværende jobb"
$ws.Range("A16").Value = "Jeg er godt fornøyd med oppgaver og ansvar i min nåværende jobb"

# A17:A18 - "forberedt_oppgaver" question text, trimmed the same way
$ws.Range("A17").Value = "Utdanningen forberedte meg godt for oppgaver og ansvar i min nåværende jobb"
$ws.Range("A18").Value = "Utdanningen forberedte meg godt for oppgaver og ansvar i min nåværende jobb"

# A19:A20 - "kompetanse_tverrprofesjonelt" question text, trimmed the same way
$ws.Range("A19").Value = "Utdanningen ga meg kompetanse i å samarbeide med andre yrkesgrupper (tverrprofesjonelt samarbeid)"
$ws.Range("A20").Value = "Utdanningen ga meg kompetanse i å samarbeide med andre yrkesgrupper (tverrprofesjonelt samarbeid)"

# A21:A22 - "valgt_samme_utdanning" question text, rephrased to a direct question
$ws.Range("A21").Value = "Hvor sannsynlig er det at du ville ha valgt samme type utdanning?"
$ws.Range("A22").Value = "Hvor sannsynlig er det at du ville ha valgt samme type utdanning?"

# A23:A24 - "valgt_samme_institusjon" question text, rephrased to a direct question
$ws.Range("A23").Value = "Hvor sannsynlig er det at du ville ha valgt samme lærested?"
$ws.Range("A24").Value = "Hvor sannsynlig er det at du ville ha valgt samme lærested?"

# --- Tidy up the view: put selection back at the top, scrolled to A1 ---
$ws.Range("A2").Select()
